$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.414.81"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.834.12"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.86"
$ws.Range("E5").Value = "  -1.64%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5281"
$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2829"
$ws.Range("E8").Value = "  -11.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06924"
$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.857.46"
$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.40"
$ws.Range("E11").Value = "  -12.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7011"
$ws.Range("E12").Value = "  -10.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07109"
$ws.Range("E13").Value = "  -7.93%  "

$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.877"
$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.25"
$ws.Range("E18").Value = "  -4.06%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007375"
$ws.Range("E19").Value = "  -6.80%  "

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.461.57"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.086.71"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.519"
$ws.Range("E22").Value = "  -1.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.850"
$ws.Range("E23").Value = "  -1.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.969"
$ws.Range("E24").Value = "  -4.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.72"
$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.677"
$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.046"
$ws.Range("E27").Value = "  -5.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.61"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.28"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.090"
$ws.Range("E30").Value = "  -1.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08738"
$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.856"
$ws.Range("E32").Value = "  -4.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04692"
$ws.Range("E33").Value = "  -3.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.901"
$ws.Range("E34").Value = "  +1.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.108"
$ws.Range("E35").Value = "  -1.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7022"
$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.072"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.171"
$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01647"
$ws.Range("E39").Value = "  -6.52%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4504"
$ws.Range("E40").Value = "  -4.75%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8730"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.41"
$ws.Range("E42").Value = "  -4.51%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.748"
$ws.Range("E44").Value = "  -2.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.973"
$ws.Range("E45").Value = "  -8.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.679"
$ws.Range("E46").Value = "  -2.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1189"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "58.88"
$ws.Range("E48").Value = "  -0.91%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.48"
$ws.Range("E49").Value = "  -3.60%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05587"
$ws.Range("E50").Value = "  -4.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8635"
$ws.Range("E51").Value = "  -2.76%  "
